$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = "RESERVING_RATE"
$ws.Range("J2").Value = "S"
$ws.Range("K2").Value = 0.04
$ws.Range("K3").Value = 0.03

$ws.Range("K1").EntireColumn.AutoFit()

$ws.Range("J3").Select()
